# edit.ps1
# Applies the "small change in read me" edit:
#   - Bumps the applied-load input cells (column H) in the "Euramet" sheet
#     calibration tables from 1 N to 500 N for both the Q3 (rows 7-21) and
#     Q1 (rows 29-43) blocks.
#   - Refreshes a handful of column F "deflection" readings that were
#     swapped/updated alongside the load change.
#   - Fills in the previously-blank rows 17-21 / 39-43 (the 0.8 / 1.0 /
#     "SetPnt" rows) with the same hm/E/F/G/H values used by their peers.
#   - Mirrors the 500 N value on the "Istruzioni Uso" sheet (B13, the scale
#     reference cell) and clears the placeholder "-" client cell (B65).

$wb  = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Euramet")

# ------------------------------------------------------------------
# "Euramet" sheet
# ------------------------------------------------------------------
    # Q3 - Negativo side (rows 7-21): calibration load-cell readings
    $ws1.Range("H7").Value = 500

    $ws1.Range("F8").Value = -157.0411782090434
    $ws1.Range("H8").Value = 500

    $ws1.Range("F9").Value = -156.5579745837849
    $ws1.Range("H9").Value = 500

    $ws1.Range("H10").Value = 500

    $ws1.Range("F11").Value = -157.0411782090434
    $ws1.Range("H11").Value = 500

    $ws1.Range("F12").Value = -156.5579745837849
    $ws1.Range("H12").Value = 500

    $ws1.Range("F13").Value = -157.0411782090434
    $ws1.Range("H13").Value = 500

    $ws1.Range("F14").Value = -157.0411782090434
    $ws1.Range("H14").Value = 500

    $ws1.Range("F15").Value = -157.0411782090434
    $ws1.Range("H15").Value = 500

    $ws1.Range("F16").Value = -157.0411782090434
    $ws1.Range("H16").Value = 500

    $ws1.Range("D17").Value = 837
    $ws1.Range("E17").Value = 0
    $ws1.Range("F17").Value = -157.0411782090434
    $ws1.Range("G17").Value = -1.1772
    $ws1.Range("H17").Value = 500

    $ws1.Range("D18").Value = 837
    $ws1.Range("E18").Value = 0
    $ws1.Range("F18").Value = -157.0411782090434
    $ws1.Range("G18").Value = -1.1772
    $ws1.Range("H18").Value = 500

    $ws1.Range("D19").Value = 837
    $ws1.Range("E19").Value = 0
    $ws1.Range("F19").Value = -157.0411782090434
    $ws1.Range("G19").Value = -1.1772
    $ws1.Range("H19").Value = 500

    $ws1.Range("D20").Value = 837
    $ws1.Range("E20").Value = 0
    $ws1.Range("F20").Value = -157.0411782090434
    $ws1.Range("G20").Value = -1.1772
    $ws1.Range("H20").Value = 500

    $ws1.Range("D21").Value = 837
    $ws1.Range("E21").Value = 0
    $ws1.Range("F21").Value = -157.0411782090434
    $ws1.Range("G21").Value = -1.1772
    $ws1.Range("H21").Value = 500

    # Q1 - Positivo side (rows 29-43): calibration load-cell readings
    $ws1.Range("H29").Value = 500

    $ws1.Range("F30").Value = -157.0411782090434
    $ws1.Range("H30").Value = 500

    $ws1.Range("H31").Value = 500

    $ws1.Range("H32").Value = 500

    $ws1.Range("H33").Value = 500

    $ws1.Range("F34").Value = -157.0411782090434
    $ws1.Range("H34").Value = 500

    $ws1.Range("F35").Value = -157.0411782090434
    $ws1.Range("H35").Value = 500

    $ws1.Range("H36").Value = 500

    $ws1.Range("F37").Value = -157.0411782090434
    $ws1.Range("H37").Value = 500

    $ws1.Range("F38").Value = -156.5579745837849
    $ws1.Range("H38").Value = 500

    $ws1.Range("D39").Value = 837
    $ws1.Range("E39").Value = 0
    $ws1.Range("F39").Value = -157.0411782090434
    $ws1.Range("G39").Value = -1.1772
    $ws1.Range("H39").Value = 500

    $ws1.Range("D40").Value = 837
    $ws1.Range("E40").Value = 0
    $ws1.Range("F40").Value = -157.0411782090434
    $ws1.Range("G40").Value = -1.1772
    $ws1.Range("H40").Value = 500

    $ws1.Range("D41").Value = 837
    $ws1.Range("E41").Value = 0
    $ws1.Range("F41").Value = -157.0411782090434
    $ws1.Range("G41").Value = -1.1772
    $ws1.Range("H41").Value = 500

    $ws1.Range("D42").Value = 837
    $ws1.Range("E42").Value = 0
    $ws1.Range("F42").Value = -157.0411782090434
    $ws1.Range("G42").Value = -1.1772
    $ws1.Range("H42").Value = 500

    $ws1.Range("D43").Value = 837
    $ws1.Range("E43").Value = 0
    $ws1.Range("F43").Value = -156.5579745837849
    $ws1.Range("G43").Value = -1.1772
    $ws1.Range("H43").Value = 500

# ------------------------------------------------------------------
# "Istruzioni Uso" sheet
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Istruzioni Uso")

# Scale reference value, mirrors the 1 -> 500 change above
$ws2.Range("B13").Value = 500

# Clear the placeholder "-" in the Cliente block
$ws2.Range("B65").ClearContents()
